# Apply "Generate Report for Handoff" update:
#  - Two new source files handed off (07b190c3-... and 94ebf8c1-...)
#  - Overview sheet gains two new rows (status "Ready for handoff"), the
#    ".localization-config" summary row moves down to make room.
#  - zh-cn / de-de detail sheets gain matching rows with handoff file,
#    handoff datetime, target file and status ("Include") details, and the
#    ".localization-config" ("Ignored") row moves down to make room.

$wb = $excel.ActiveWorkbook
$missing = [System.Reflection.Missing]::Value

$mdBase  = "https://github.com/OpenLocalizationTest/oltest/blob/8779e9ee3d0bffabe848effc99dbc3fd0d1e44d8"
$zhBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/912ad3e9aaf04aa90b1a86126d3b0430aa1d043b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang"
$deBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4753edc27d6359e4b99aa68b2e93e878f71fb273/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang"

$file1 = "07b190c3-b3e5-4230-a29e-3ad6e0bd6d5e.md"
$file1ZhXlf = "07b190c3-b3e5-4230-a29e-3ad6e0bd6d5e.be54133101a2513053ae1c0ff0a5a91473344b93.zh-cn.xlf"
$file1DeXlf = "07b190c3-b3e5-4230-a29e-3ad6e0bd6d5e.be54133101a2513053ae1c0ff0a5a91473344b93.de-de.xlf"

$file2 = "94ebf8c1-2053-4368-b53c-af4bd15eb411.md"
$file2ZhXlf = "94ebf8c1-2053-4368-b53c-af4bd15eb411.6914bce98cf6cb0bc2e095cd9253fd8559c610bc.zh-cn.xlf"
$file2DeXlf = "94ebf8c1-2053-4368-b53c-af4bd15eb411.6914bce98cf6cb0bc2e095cd9253fd8559c610bc.de-de.xlf"

$zhHandoffDateTime = "2016-01-27 02:11:28"
$deHandoffDateTime = "2016-01-27 02:11:39"
$epoch = "0001-01-01 00:00:00"

# -----------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A1").Hyperlinks.Delete()

$ws.Range("A4").Value = $file1
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"

$ws.Range("A5").Value = $file2
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"

$ws.Range("A6").Value = ".localization-config"
$ws.Range("B6").Value = "Not to be localized"
$ws.Range("C6").Value = "Not to be localized"

$ws.Hyperlinks.Add($ws.Range("A2"), "$mdBase/e2e/8c38aaeb-002d-4fac-a6f2-1a8422a9446b.md", $missing, $missing, "8c38aaeb-002d-4fac-a6f2-1a8422a9446b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "$mdBase/e2e/e3889229-4ba8-4197-a20e-b3f3b762a011.md", $missing, $missing, "e3889229-4ba8-4197-a20e-b3f3b762a011.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "$mdBase/e2e/$file1", $missing, $missing, $file1) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "$mdBase/e2e/$file2", $missing, $missing, $file2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "$mdBase/.localization-config", $missing, $missing, ".localization-config") | Out-Null

# -----------------------------------------------------------------
# Sheet "zh-cn" detail table
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A1").Hyperlinks.Delete()

$ws.Range("A4").Value = $file1
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = $file1ZhXlf
$ws.Range("D4").Value = $zhHandoffDateTime
$ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G4").Value = $epoch
$ws.Range("H4").Value = "Include"

$ws.Range("A5").Value = $file2
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = $file2ZhXlf
$ws.Range("D5").Value = $zhHandoffDateTime
$ws.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G5").Value = $epoch
$ws.Range("H5").Value = "Include"

$ws.Range("A6").Value = ".localization-config"
$ws.Range("B6").Value = "Not to be localized"
$ws.Range("D6").Value = $epoch
$ws.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G6").Value = $epoch
$ws.Range("H6").Value = "Ignored"

$ws.Hyperlinks.Add($ws.Range("A2"), "$mdBase/e2e/8c38aaeb-002d-4fac-a6f2-1a8422a9446b.md", $missing, $missing, "8c38aaeb-002d-4fac-a6f2-1a8422a9446b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "$zhBase/8c38aaeb-002d-4fac-a6f2-1a8422a9446b.da2348ca289544631f21abdb080fb33482517d29.zh-cn.xlf", $missing, $missing, "8c38aaeb-002d-4fac-a6f2-1a8422a9446b.da2348ca289544631f21abdb080fb33482517d29.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "$mdBase/e2e/e3889229-4ba8-4197-a20e-b3f3b762a011.md", $missing, $missing, "e3889229-4ba8-4197-a20e-b3f3b762a011.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "$zhBase/e3889229-4ba8-4197-a20e-b3f3b762a011.291e4de995e596ac7a17757c820cf777ce282ae4.zh-cn.xlf", $missing, $missing, "e3889229-4ba8-4197-a20e-b3f3b762a011.291e4de995e596ac7a17757c820cf777ce282ae4.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "$mdBase/e2e/$file1", $missing, $missing, $file1) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "$zhBase/$file1ZhXlf", $missing, $missing, $file1ZhXlf) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "$mdBase/e2e/$file2", $missing, $missing, $file2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "$zhBase/$file2ZhXlf", $missing, $missing, $file2ZhXlf) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "$mdBase/.localization-config", $missing, $missing, ".localization-config") | Out-Null

# -----------------------------------------------------------------
# Sheet "de-de" detail table
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A1").Hyperlinks.Delete()

$ws.Range("A4").Value = $file1
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = $file1DeXlf
$ws.Range("D4").Value = $deHandoffDateTime
$ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G4").Value = $epoch
$ws.Range("H4").Value = "Include"

$ws.Range("A5").Value = $file2
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = $file2DeXlf
$ws.Range("D5").Value = $deHandoffDateTime
$ws.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G5").Value = $epoch
$ws.Range("H5").Value = "Include"

$ws.Range("A6").Value = ".localization-config"
$ws.Range("B6").Value = "Not to be localized"
$ws.Range("D6").Value = $epoch
$ws.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G6").Value = $epoch
$ws.Range("H6").Value = "Ignored"

$ws.Hyperlinks.Add($ws.Range("A2"), "$mdBase/e2e/8c38aaeb-002d-4fac-a6f2-1a8422a9446b.md", $missing, $missing, "8c38aaeb-002d-4fac-a6f2-1a8422a9446b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "$deBase/8c38aaeb-002d-4fac-a6f2-1a8422a9446b.da2348ca289544631f21abdb080fb33482517d29.de-de.xlf", $missing, $missing, "8c38aaeb-002d-4fac-a6f2-1a8422a9446b.da2348ca289544631f21abdb080fb33482517d29.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "$mdBase/e2e/e3889229-4ba8-4197-a20e-b3f3b762a011.md", $missing, $missing, "e3889229-4ba8-4197-a20e-b3f3b762a011.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "$deBase/e3889229-4ba8-4197-a20e-b3f3b762a011.291e4de995e596ac7a17757c820cf777ce282ae4.de-de.xlf", $missing, $missing, "e3889229-4ba8-4197-a20e-b3f3b762a011.291e4de995e596ac7a17757c820cf777ce282ae4.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "$mdBase/e2e/$file1", $missing, $missing, $file1) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "$deBase/$file1DeXlf", $missing, $missing, $file1DeXlf) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "$mdBase/e2e/$file2", $missing, $missing, $file2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "$deBase/$file2DeXlf", $missing, $missing, $file2DeXlf) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "$mdBase/.localization-config", $missing, $missing, ".localization-config") | Out-Null

Write-Host "Report generated for handoff."
